$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 4431
$ws.Range('K3').Value = 4543
$ws.Range('C4').Value = 1850
$ws.Range('K4').Value = 914
$ws.Range('K6').Value = 5121
$ws.Range('C7').Value = 28395
$ws.Range('K7').Value = 15337

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 133
$ws.Range('K7').Value = 449
$ws.Range('K8').Value = 1024
$ws.Range('K11').Value = 301
$ws.Range('K14').Value = 86
$ws.Range('K15').Value = 155
$ws.Range('K19').Value = 465
$ws.Range('K20').Value = 351
$ws.Range('K29').Value = 808
$ws.Range('K30').Value = 55
$ws.Range('K31').Value = 167
$ws.Range('K33').Value = 641
$ws.Range('K34').Value = 79
$ws.Range('K36').Value = 193
$ws.Range('K37').Value = 521
$ws.Range('K42').Value = 568
$ws.Range('K43').Value = 138
$ws.Range('K44').Value = 136
$ws.Range('K47').Value = 96
$ws.Range('K48').Value = 199
$ws.Range('K49').Value = 88
$ws.Range('K51').Value = 195
$ws.Range('K52').Value = 406
$ws.Range('K53').Value = 204
$ws.Range('K54').Value = 287
$ws.Range('K60').Value = 99
$ws.Range('C63').Value = 279
$ws.Range('K63').Value = 48
$ws.Range('K64').Value = 94
$ws.Range('K65').Value = 347
$ws.Range('K67').Value = 592
$ws.Range('K68').Value = 39
$ws.Range('K72').Value = 70
$ws.Range('K75').Value = 53
$ws.Range('K76').Value = 212
$ws.Range('K79').Value = 387
$ws.Range('K83').Value = 328
$ws.Range('K85').Value = 688
$ws.Range('K88').Value = 178
$ws.Range('K89').Value = 221
$ws.Range('K90').Value = 142
$ws.Range('K93').Value = 57
$ws.Range('K95').Value = 273
$ws.Range('K96').Value = 171
$ws.Range('K98').Value = 77
$ws.Range('C101').Value = 28395
$ws.Range('K101').Value = 15337

# Sheet 3: Bridgeport
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('K6').Value = 31
$ws.Range('K7').Value = 86

# Sheet 4: West Ridge
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K4').Value = 7
$ws.Range('K7').Value = 171

# Sheet 5: Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 157
$ws.Range('K3').Value = 142
$ws.Range('K7').Value = 449

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K2').Value = 99
$ws.Range('K7').Value = 301

# Sheet 7: Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K3').Value = 69
$ws.Range('K7').Value = 221

# Sheet 8: South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 241
$ws.Range('K3').Value = 230
$ws.Range('K7').Value = 688

# Sheet 9: Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K6').Value = 153
$ws.Range('K7').Value = 406

# Sheet 11: Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K6').Value = 94
$ws.Range('K7').Value = 204

# Sheet 12: Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K3').Value = 307
$ws.Range('K4').Value = 59
$ws.Range('K6').Value = 345
$ws.Range('K7').Value = 1024

# Sheet 13: South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 113
$ws.Range('K3').Value = 119
$ws.Range('K7').Value = 328

# Sheet 14: Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K3').Value = 238
$ws.Range('K6').Value = 186
$ws.Range('K7').Value = 641

# Sheet 15: West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K2').Value = 91
$ws.Range('K3').Value = 92
$ws.Range('K7').Value = 273

# Sheet 16: Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 145
$ws.Range('K3').Value = 170
$ws.Range('K6').Value = 157
$ws.Range('K7').Value = 521

# Sheet 17: New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range('K3').Value = 86
$ws.Range('K6').Value = 137
$ws.Range('K7').Value = 347

# Sheet 19: Fuller Park
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('K3').Value = 21
$ws.Range('K7').Value = 55

# Sheet 20: Gage Park
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K2').Value = 58
$ws.Range('K7').Value = 167

# Sheet 21: North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 167
$ws.Range('K4').Value = 33
$ws.Range('K7').Value = 592

# Sheet 23: Lincoln Park
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K6').Value = 50
$ws.Range('K7').Value = 88

# Sheet 24: Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K6').Value = 146
$ws.Range('K7').Value = 287

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K3').Value = 287
$ws.Range('K6').Value = 226
$ws.Range('K7').Value = 808

# Sheet 26: Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K6').Value = 98
$ws.Range('K7').Value = 199

# Sheet 27: Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 144
$ws.Range('K3').Value = 142
$ws.Range('K6').Value = 146
$ws.Range('K7').Value = 465

# Sheet 28: Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K2').Value = 30
$ws.Range('K7').Value = 136

# Sheet 29: River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range('K2').Value = 43
$ws.Range('K3').Value = 40
$ws.Range('K6').Value = 117
$ws.Range('K7').Value = 212

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 153
$ws.Range('K3').Value = 175
$ws.Range('K4').Value = 23
$ws.Range('K6').Value = 213
$ws.Range('K7').Value = 568

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K3').Value = 127
$ws.Range('K7').Value = 387

# Sheet 43: Near South Side
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K3').Value = 29
$ws.Range('K7').Value = 94

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 121
$ws.Range('K6').Value = 103
$ws.Range('K7').Value = 351

# Sheet 47: Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K3').Value = 54
$ws.Range('K7').Value = 193

# Sheet 48: West Lawn
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('K2').Value = 20
$ws.Range('K7').Value = 57

# Sheet 50: Garfield Ridge
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K4').Value = 4
$ws.Range('K7').Value = 79

# Sheet 53: Kenwood
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K2').Value = 31
$ws.Range('K7').Value = 96

# Sheet 54: Brighton Park
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K4').Value = 11
$ws.Range('K6').Value = 47
$ws.Range('K7').Value = 155

# Sheet 55: Wicker Park
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('K2').Value = 11
$ws.Range('K7').Value = 77

# Sheet 64: Albany Park
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K2').Value = 43
$ws.Range('K7').Value = 133

# Sheet 68: United Center
$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K2').Value = 42
$ws.Range('K7').Value = 178

# Sheet 73: Pullman
$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('K6').Value = 8
$ws.Range('K7').Value = 53

# Sheet 74: Washington Heights
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K3').Value = 46
$ws.Range('K7').Value = 142

# Sheet 75: Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K2').Value = 53
$ws.Range('K7').Value = 195

# Sheet 76: North Park
$ws = $wb.Worksheets.Item('North Park')
$ws.Range('K6').Value = 10
$ws.Range('K7').Value = 39

# Sheet 78: Morgan Park
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('K2').Value = 36
$ws.Range('K7').Value = 99

# Sheet 79: Hyde Park
$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('K4').Value = 17
$ws.Range('K7').Value = 138

# Sheet 82: Old Town
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('K2').Value = 9
$ws.Range('K6').Value = 39
$ws.Range('K7').Value = 70
